$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.8900471604056629
$ws.Range("J2").Value = 0.8900471604056629
$ws.Range("S2").Value = 0.8900471604056629
$ws.Range("T2").Value = 0.8900471604056629

# Row 3 updates
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.007903666666666666
$ws.Range("H3").Value = 0.023711
$ws.Range("I3").Value = 0.109952839594337
$ws.Range("J3").Value = 0.109952839594337
$ws.Range("Q3").Value = 0.0002133489434444444
$ws.Range("R3").Value = 0.001920140491
$ws.Range("S3").Value = 0.109952839594337
$ws.Range("T3").Value = 0.109952839594337
